# Add a new paragraph for the Liberation Mono font right after the
# existing "Liberation Sans" paragraph, per the commit:
#   "Added Liberation Mono to fonts needed"

$d = $word.ActiveDocument

# Locate the paragraph that announces the Liberation Sans font download
# (robust to exact paragraph index even if the document shape changes).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Download the Liberation Sans font*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a brand-new paragraph directly after it, then fill it in --
    # this keeps the original paragraph/run completely untouched and
    # matches the new <w:p><w:r><w:t>.../w:t></w:r></w:p> block in the diff.
    $target.Range.InsertParagraphAfter()
    $newRange = $target.Next().Range
    $newRange.Text = "Download the Liberation Mono font from http://www.fontsquirrel.com/fonts/Liberation-Mono and install."
}
